# Generate Report for Handoff
# Refresh the localization-status workbook with a new handoff run:
#  - new GUID for the source/target file basenames
#    (73e91299-47c5-4767-9fda-a23aa462bcfa -> fc95439c-1e7a-4315-905b-f0bd0defa416)
#  - new content hash for the generated .xlf targets
#    (ea2495e14a87efca83c9f4352d5274203dd94048 -> 0ac318e34844fed6e5954a1deb49943210cee43f)
#  - refreshed handoff timestamps

$wb = $excel.ActiveWorkbook

$oldGuid = "73e91299-47c5-4767-9fda-a23aa462bcfa"
$newGuid = "fc95439c-1e7a-4315-905b-f0bd0defa416"
$oldHash = "ea2495e14a87efca83c9f4352d5274203dd94048"
$newHash = "0ac318e34844fed6e5954a1deb49943210cee43f"

$newMdName    = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHash.de-de.xlf"

$newHandoffDateBad = "2016-44-18 08:44:26"   # Sheet "Overview" Latest Handoff Date (mirrors source malformed value)
$newZhDateTime      = "2016-03-18 08:44:23"  # Sheet "zh-cn" Latest Handoff Datetime
$newDeDateTime      = "2016-03-18 08:44:26"  # Sheet "de-de" Latest Handoff Datetime

# Hyperlink targets are untouched by this edit -- only the visible display
# text / backing shared-string changes. The runtime's Hyperlinks collection
# can't be updated in place (and Address/TextToDisplay read back empty), so
# rebuild each sheet's hyperlinks (delete-all then re-add in the same order)
# using the original target URLs, pointing at the refreshed display text.

$mdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/0180e6a828b587be12b3b71359669a6bfbdbf525/e2e/$oldGuid.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a790f4e23b7ce53336d5e104066d2e2a070228b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3853c207a4ca39fc26a615f43037656b5c5fa96b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf"

# ---- Sheet 1: Overview ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", $newMdName) | Out-Null
$wsOverview.Range("D2").Value = $newHandoffDateBad

# ---- Sheet 2: zh-cn ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl, "", "", $newMdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $mdUrl, "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfUrl, "", "", $newZhXlfName) | Out-Null
$wsZh.Range("E2").Value = $newZhDateTime

# ---- Sheet 3: de-de ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl, "", "", $newMdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $mdUrl, "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfUrl, "", "", $newDeXlfName) | Out-Null
$wsDe.Range("E2").Value = $newDeDateTime
